$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 3999.875
$ws.Range("I82").Value = 2666.5
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 7999.5
$ws.Range("L82").Value = 24000
$ws.Range("M82").Value = -7593.5
$ws.Range("N82").Value = -24812
$ws.Range("H85").Value = 3999.875
$ws.Range("I85").Value = 2666.5
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 7999.5
$ws.Range("L85").Value = 24000
$ws.Range("M85").Value = -6595.5
$ws.Range("N85").Value = -26808
$ws.Range("H121").Value = 1236
$ws.Range("J121").Value = 1641
$ws.Range("L121").Value = 4923
$ws.Range("N121").Value = -8417
$ws.Range("H137").Value = 5004800
$ws.Range("I137").Value = 7698084.5
$ws.Range("K137").Value = 23094253.5
$ws.Range("M137").Value = -23091703.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 142.14285
$ws.Range("J4").Value = 174.25
$ws.Range("L4").Value = 174.25
$ws.Range("N4").Value = -406.25
$ws.Range("H32").Value = 6216.8423
$ws.Range("I32").Value = 3699.2554
$ws.Range("J32").Value = 18049.5
$ws.Range("K32").Value = 3699.2554
$ws.Range("L32").Value = 18049.5
$ws.Range("M32").Value = -3412.2554
$ws.Range("N32").Value = -18623.5
$ws.Range("H61").Value = 2699.0908
$ws.Range("I61").Value = 1632.2222
$ws.Range("J61").Value = 7500
$ws.Range("K61").Value = 1632.2222
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -1420.2222
$ws.Range("N61").Value = -7924
$ws.Range("H132").Value = 30307672
$ws.Range("I132").Value = 47623484
$ws.Range("J132").Value = 5002.1665
$ws.Range("K132").Value = 142870452
$ws.Range("L132").Value = 15006.4995
$ws.Range("M132").Value = -142867922
$ws.Range("N132").Value = -20066.4995
$ws.Range("H136").Value = 2699.0908
$ws.Range("I136").Value = 1632.2222
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 4896.6666
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -2346.6666
$ws.Range("N136").Value = -27600

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1585.9584
$ws.Range("I105").Value = 1314.7059
$ws.Range("K105").Value = 1314.7059
$ws.Range("M105").Value = 432.2941000000001
$ws.Range("H134").Value = 2994.0952
$ws.Range("I134").Value = 1773.5
$ws.Range("J134").Value = 6900
$ws.Range("K134").Value = 5320.5
$ws.Range("L134").Value = 20700
$ws.Range("M134").Value = -2785.5
$ws.Range("N134").Value = -25770

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 20002780
$ws.Range("I58").Value = 1367
$ws.Range("J58").Value = 50004900
$ws.Range("K58").Value = 1367
$ws.Range("L58").Value = 50004900
$ws.Range("M58").Value = -1164
$ws.Range("N58").Value = -50005306
$ws.Range("H132").Value = 3136.698
$ws.Range("I132").Value = 2074.7666
$ws.Range("J132").Value = 4521.826
$ws.Range("K132").Value = 6224.2998
$ws.Range("L132").Value = 13565.478
$ws.Range("M132").Value = -3694.2998
$ws.Range("N132").Value = -18625.478
$ws.Range("H136").Value = 20002780
$ws.Range("I136").Value = 1367
$ws.Range("J136").Value = 50004900
$ws.Range("K136").Value = 4101
$ws.Range("L136").Value = 150014700
$ws.Range("M136").Value = -1551
$ws.Range("N136").Value = -150019800
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2149.2666
$ws.Range("I137").Value = 1423.9
$ws.Range("J137").Value = 3600
$ws.Range("K137").Value = 4271.700000000001
$ws.Range("L137").Value = 10800
$ws.Range("M137").Value = 828.2999999999993
$ws.Range("N137").Value = -21000

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4609.9
$ws.Range("I122").Value = 6056.857
$ws.Range("K122").Value = 18170.571
$ws.Range("M122").Value = -15720.571
$ws.Range("H132").Value = 3049.7026
$ws.Range("I132").Value = 2471.348
$ws.Range("J132").Value = 3999.8572
$ws.Range("K132").Value = 7414.044
$ws.Range("L132").Value = 11999.5716
$ws.Range("M132").Value = -4884.044
$ws.Range("N132").Value = -17059.5716

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("H16").Value = 673.5
$ws.Range("I16").Value = 732.1667
$ws.Range("K16").Value = 732.1667
$ws.Range("M16").Value = -562.1667
$ws.Range("H22").Value = 55556810
$ws.Range("J22").Value = 1418.3334
$ws.Range("L22").Value = 1418.3334
$ws.Range("N22").Value = -2008.3334
$ws.Range("H27").Value = 55556810
$ws.Range("J27").Value = 1418.3334
$ws.Range("L27").Value = 1418.3334
$ws.Range("N27").Value = -1632.3334
$ws.Range("H46").Value = 3224.8333
$ws.Range("I46").Value = 702
$ws.Range("J46").Value = 5026.857
$ws.Range("K46").Value = 702
$ws.Range("L46").Value = 5026.857
$ws.Range("M46").Value = -514
$ws.Range("N46").Value = -5402.857
$ws.Range("H82").Value = 2208.1052
$ws.Range("I82").Value = 1875.1111
$ws.Range("J82").Value = 2507.8
$ws.Range("K82").Value = 1875.1111
$ws.Range("L82").Value = 2507.8
$ws.Range("M82").Value = -1514.1111
$ws.Range("N82").Value = -3229.8
$ws.Range("H85").Value = 2208.1052
$ws.Range("I85").Value = 1875.1111
$ws.Range("J85").Value = 2507.8
$ws.Range("K85").Value = 1875.1111
$ws.Range("L85").Value = 2507.8
$ws.Range("M85").Value = -627.1111000000001
$ws.Range("N85").Value = -5003.8
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H132").Value = 2432.238
$ws.Range("I132").Value = 1331.7037
$ws.Range("J132").Value = 4413.2
$ws.Range("K132").Value = 3995.1111
$ws.Range("L132").Value = 13239.6
$ws.Range("M132").Value = -1465.1111
$ws.Range("N132").Value = -18299.6
$ws.Range("H136").Value = 4350495.5
$ws.Range("I136").Value = 7694169
$ws.Range("J136").Value = 3719.5
$ws.Range("K136").Value = 23082507
$ws.Range("L136").Value = 11158.5
$ws.Range("M136").Value = -23079957
$ws.Range("N136").Value = -16258.5
$ws.Range("M7").ClearContents()
$ws.Range("M126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 28233.334
$ws.Range("I82").Value = 15000
$ws.Range("J82").Value = 29887.5
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 29887.5
$ws.Range("M82").Value = -14617
$ws.Range("N82").Value = -30653.5
$ws.Range("H85").Value = 28233.334
$ws.Range("I85").Value = 15000
$ws.Range("J85").Value = 29887.5
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 29887.5
$ws.Range("M85").Value = -13674
$ws.Range("N85").Value = -32539.5
$ws.Range("H122").Value = 346562.28
$ws.Range("I122").Value = 418155.12
$ws.Range("J122").Value = 2916.6
$ws.Range("K122").Value = 1254465.36
$ws.Range("L122").Value = 8749.799999999999
$ws.Range("M122").Value = -1252015.36
$ws.Range("N122").Value = -13649.8
$ws.Range("H132").Value = 288073.22
$ws.Range("I132").Value = 439661.12
$ws.Range("K132").Value = 1318983.36
$ws.Range("M132").Value = -1316453.36
$ws.Range("H136").Value = 2207.875
$ws.Range("I136").Value = 1443.8334
$ws.Range("J136").Value = 4500
$ws.Range("K136").Value = 4331.5002
$ws.Range("L136").Value = 13500
$ws.Range("M136").Value = -1781.5002
$ws.Range("N136").Value = -18600
